$p = $ppt.ActivePresentation

# The deck's single slide master (and every slide, since they all share it)
# is wired to ppt/theme/theme2.xml, which currently holds the "Integral" /
# "Red Violet" theme. The target edit swaps the Integral theme into
# ppt/theme/theme1.xml (used only by the notes master) and puts the
# default "Office Theme" / "Office" color palette into ppt/theme/theme2.xml
# -- i.e. the presentation's visible theme reverts to the stock Office
# colors. We reproduce that visible/semantic effect by rewriting the 12
# theme colors of the active theme (reached via a slide's
# ThemeColorScheme, which maps 1:1 onto dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink) to the standard Office palette.

$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}
